$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set Price column (D) to Text format first so numeric-looking values
# (e.g. "48.40", "1.000") are preserved as text rather than being
# converted to numbers, matching the original inline-string cells.
$ws.Range("D2:D51").NumberFormat = "@"

# Update Price (D) and Volume(1h) (E) values per row
$ws.Range("D2").Value = "26.559.23"
$ws.Range("E2").Value = "  +6.96%  "
$ws.Range("D3").Value = "1.726.05"
$ws.Range("E3").Value = "  +3.63%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("E5").Value = "  +0.92%  "
$ws.Range("D6").Value = "0.9994"
$ws.Range("E6").Value = "  -0.27%  "
$ws.Range("D7").Value = "0.3715"
$ws.Range("E7").Value = "  +1.89%  "
$ws.Range("D8").Value = "48.40"
$ws.Range("E8").Value = "  +2.41%  "
$ws.Range("D9").Value = "0.3361"
$ws.Range("E9").Value = "  +3.08%  "
$ws.Range("E10").Value = "  +4.08%  "
$ws.Range("D11").Value = "0.07407"
$ws.Range("E11").Value = "  +4.64%  "
$ws.Range("D12").Value = "1.000"
$ws.Range("E12").Value = "  -0.04%  "
$ws.Range("D13").Value = "6.379"
$ws.Range("E13").Value = "  +5.03%  "
$ws.Range("D14").Value = "20.09"
$ws.Range("E14").Value = "  +3.06%  "
$ws.Range("D15").Value = "7.049"
$ws.Range("E15").Value = "  +6.81%  "
$ws.Range("D16").Value = "1.727.32"
$ws.Range("E16").Value = "  +3.73%  "
$ws.Range("E17").Value = "  +2.05%  "
$ws.Range("D18").Value = "0.06627"
$ws.Range("E18").Value = "  -0.27%  "
$ws.Range("D19").Value = "81.91"
$ws.Range("E19").Value = "  +4.35%  "
$ws.Range("D20").Value = "0.9986"
$ws.Range("E20").Value = "  -0.28%  "
$ws.Range("D21").Value = "16.55"
$ws.Range("E21").Value = "  +4.79%  "
$ws.Range("D22").Value = "6.133"
$ws.Range("E22").Value = "  +3.50%  "
$ws.Range("E23").Value = "  +1.74%  "
$ws.Range("D24").Value = "26.543.27"
$ws.Range("E24").Value = "  +7.00%  "
$ws.Range("D25").Value = "2.446"
$ws.Range("E25").Value = "  -0.90%  "
$ws.Range("D26").Value = "1.417"
$ws.Range("E26").Value = "  +21.19%  "
$ws.Range("D27").Value = "2.397"
$ws.Range("E27").Value = "  -1.57%  "
$ws.Range("D28").Value = "151.41"
$ws.Range("E28").Value = "  +1.89%  "
$ws.Range("D29").Value = "19.43"
$ws.Range("E29").Value = "  +4.04%  "
$ws.Range("D30").Value = "1.917.48"
$ws.Range("E30").Value = "  +3.81%  "
$ws.Range("D31").Value = "131.31"
$ws.Range("E31").Value = "  +4.16%  "
$ws.Range("D32").Value = "4.110"
$ws.Range("E32").Value = "  +1.05%  "
$ws.Range("D33").Value = "5.965"
$ws.Range("E33").Value = "  +4.89%  "
$ws.Range("D34").Value = "0.08609"
$ws.Range("E34").Value = "  +1.76%  "
$ws.Range("D35").Value = "1.694"
$ws.Range("E35").Value = "  +2.95%  "
$ws.Range("D36").Value = "12.73"
$ws.Range("E36").Value = "  +4.66%  "
$ws.Range("D37").Value = "5.355"
$ws.Range("E37").Value = "  +3.62%  "
$ws.Range("D38").Value = "0.02330"
$ws.Range("E38").Value = "  +2.07%  "
$ws.Range("D39").Value = "0.06203"
$ws.Range("E39").Value = "  -0.67%  "
$ws.Range("D40").Value = "0.2147"
$ws.Range("E40").Value = "  +2.96%  "
$ws.Range("E41").Value = "  +2.17%  "
$ws.Range("D42").Value = "1.221"
$ws.Range("E42").Value = "  -1.78%  "
$ws.Range("D43").Value = "0.6192"
$ws.Range("E43").Value = "  +4.31%  "
$ws.Range("D44").Value = "14.24"
$ws.Range("E44").Value = "  +6.18%  "
$ws.Range("D45").Value = "0.9989"
$ws.Range("E45").Value = "  -0.19%  "
$ws.Range("D46").Value = "3.905"
$ws.Range("E46").Value = "  +1.55%  "
$ws.Range("D47").Value = "0.6010"
$ws.Range("E47").Value = "  +6.31%  "
$ws.Range("D48").Value = "128.95"
$ws.Range("E48").Value = "  +2.62%  "
$ws.Range("D49").Value = "2.040"
$ws.Range("E49").Value = "  +4.74%  "
$ws.Range("D50").Value = "0.07166"
$ws.Range("E50").Value = "  +2.81%  "
$ws.Range("D51").Value = "76.84"
$ws.Range("E51").Value = "  +2.20%  "

# Restore default (Normal) style on column D so no stray number format
# is left behind on the cells (keeps styling identical to original).
$ws.Range("D2:D51").Style = "Normal"
